# Applies the LinuxForHealth/alvearie-fhir-ig deployment update to the
# "employee-wage-basis" StructureDefinition workbook.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-wage-basis"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates --------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 (the "Extension" element row) no longer carries the ele-1/ext-1
# constraint text in the Constraint(s) column (AI) -- it has been relocated
# to the "Extension.extension" row (row 4) in this revision.
$elements.Range("AI2").Value = ""
